# Scheduled-runner update: refresh the cached marketboard-derived figures
# (currentAveragePrice / *NQ / *HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns
# H-N) on a handful of rows across the per-job Leve tables. Values are
# re-synced one row at a time; a few rows also gain/lose a trailing
# LeveProfit cell (H/J/K/L etc. always present, M/N only written when
# profit is computable for that side).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2848.3774
$ws.Cells.Item(17, 10).Value = 2848.3774
$ws.Cells.Item(17, 12).Value = 8545.1322
$ws.Cells.Item(17, 14).Value = -8881.1322

$ws.Cells.Item(40, 8).Value = 2056.6086
$ws.Cells.Item(40, 9).Value = 2227.889
$ws.Cells.Item(40, 10).Value = 1440
$ws.Cells.Item(40, 11).Value = 2227.889
$ws.Cells.Item(40, 12).Value = 1440
$ws.Cells.Item(40, 13).Value = -2052.889
$ws.Cells.Item(40, 14).Value = -1790

$ws.Cells.Item(100, 8).Value = 1441.5333
$ws.Cells.Item(100, 9).Value = 1415.4
$ws.Cells.Item(100, 10).Value = 1493.8
$ws.Cells.Item(100, 11).Value = 1415.4
$ws.Cells.Item(100, 12).Value = 1493.8
$ws.Cells.Item(100, 13).Value = -874.4000000000001
$ws.Cells.Item(100, 14).Value = -2575.8

$ws.Cells.Item(112, 8).Value = 1671.1786
$ws.Cells.Item(112, 9).Value = 495
$ws.Cells.Item(112, 11).Value = 1485
$ws.Cells.Item(112, 13).Value = -377

$ws.Cells.Item(132, 8).Value = 2701.9524
$ws.Cells.Item(132, 9).Value = 3029.7222
$ws.Cells.Item(132, 10).Value = 735.3333
$ws.Cells.Item(132, 11).Value = 9089.1666
$ws.Cells.Item(132, 12).Value = 2205.9999
$ws.Cells.Item(132, 13).Value = -6559.1666
$ws.Cells.Item(132, 14).Value = -7265.9999

$ws.Cells.Item(137, 8).Value = 2307.261
$ws.Cells.Item(137, 9).Value = 1607.8
$ws.Cells.Item(137, 10).Value = 3618.75
$ws.Cells.Item(137, 11).Value = 4823.4
$ws.Cells.Item(137, 12).Value = 10856.25
$ws.Cells.Item(137, 13).Value = -2273.4
$ws.Cells.Item(137, 14).Value = -15956.25

$ws.Cells.Item(138, 8).Value = 3477.5303
$ws.Cells.Item(138, 9).Value = 1971.7391
$ws.Cells.Item(138, 10).Value = 4282.9536
$ws.Cells.Item(138, 11).Value = 5915.2173
$ws.Cells.Item(138, 12).Value = 12848.8608
$ws.Cells.Item(138, 13).Value = -775.2173000000003
$ws.Cells.Item(138, 14).Value = -23128.8608

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 8901.75
$ws.Cells.Item(26, 9).Value = 3003.5
$ws.Cells.Item(26, 10).Value = 14800
$ws.Cells.Item(26, 11).Value = 3003.5
$ws.Cells.Item(26, 12).Value = 14800
$ws.Cells.Item(26, 13).Value = -2673.5
$ws.Cells.Item(26, 14).Value = -15460

$ws.Cells.Item(32, 8).Value = 18426.41
$ws.Cells.Item(32, 9).Value = 19637.635
$ws.Cells.Item(32, 10).Value = 10795.7
$ws.Cells.Item(32, 11).Value = 19637.635
$ws.Cells.Item(32, 12).Value = 10795.7
$ws.Cells.Item(32, 13).Value = -19350.635
$ws.Cells.Item(32, 14).Value = -11369.7

$ws.Cells.Item(45, 8).Value = 1710.409
$ws.Cells.Item(45, 9).Value = 1650.1875
$ws.Cells.Item(45, 10).Value = 1871
$ws.Cells.Item(45, 11).Value = 1650.1875
$ws.Cells.Item(45, 12).Value = 1871
$ws.Cells.Item(45, 13).Value = -1273.1875
$ws.Cells.Item(45, 14).Value = -2625

$ws.Cells.Item(61, 8).Value = 10383.595
$ws.Cells.Item(61, 9).Value = 8056.44
$ws.Cells.Item(61, 10).Value = 15231.833
$ws.Cells.Item(61, 11).Value = 8056.44
$ws.Cells.Item(61, 12).Value = 15231.833
$ws.Cells.Item(61, 13).Value = -7844.44
$ws.Cells.Item(61, 14).Value = -15655.833

$ws.Cells.Item(102, 8).Value = 1979.7916
$ws.Cells.Item(102, 9).Value = 1495.8
$ws.Cells.Item(102, 10).Value = 4399.75
$ws.Cells.Item(102, 11).Value = 1495.8
$ws.Cells.Item(102, 12).Value = 4399.75
$ws.Cells.Item(102, 13).Value = 126.2
$ws.Cells.Item(102, 14).Value = -7643.75

$ws.Cells.Item(110, 8).Value = 1138.1936
$ws.Cells.Item(110, 9).Value = 1099.4073
$ws.Cells.Item(110, 11).Value = 1099.4073
$ws.Cells.Item(110, 13).Value = 945.5926999999999

$ws.Cells.Item(132, 8).Value = 5273.3423
$ws.Cells.Item(132, 9).Value = 2206.3914
$ws.Cells.Item(132, 10).Value = 9976
$ws.Cells.Item(132, 11).Value = 6619.174199999999
$ws.Cells.Item(132, 12).Value = 29928
$ws.Cells.Item(132, 13).Value = -4089.174199999999
$ws.Cells.Item(132, 14).Value = -34988

$ws.Cells.Item(136, 8).Value = 10383.595
$ws.Cells.Item(136, 9).Value = 8056.44
$ws.Cells.Item(136, 10).Value = 15231.833
$ws.Cells.Item(136, 11).Value = 24169.32
$ws.Cells.Item(136, 12).Value = 45695.499
$ws.Cells.Item(136, 13).Value = -21619.32
$ws.Cells.Item(136, 14).Value = -50795.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 51476.57
$ws.Cells.Item(134, 9).Value = 4384.1665
$ws.Cells.Item(134, 10).Value = 114266.445
$ws.Cells.Item(134, 11).Value = 13152.4995
$ws.Cells.Item(134, 12).Value = 342799.335
$ws.Cells.Item(134, 13).Value = -10617.4995
$ws.Cells.Item(134, 14).Value = -347869.335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(15, 8).Value = 19800
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 19800
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 19800
$ws.Cells.Item(15, 13).ClearContents()
$ws.Cells.Item(15, 14).Value = -20140

$ws.Cells.Item(29, 8).Value = 24800
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 13).ClearContents()

$ws.Cells.Item(31, 8).Value = 5221.514
$ws.Cells.Item(31, 9).Value = 5683.5
$ws.Cells.Item(31, 10).Value = 4213.5454
$ws.Cells.Item(31, 11).Value = 5683.5
$ws.Cells.Item(31, 12).Value = 4213.5454
$ws.Cells.Item(31, 13).Value = -5388.5
$ws.Cells.Item(31, 14).Value = -4803.5454

$ws.Cells.Item(34, 8).Value = 5221.514
$ws.Cells.Item(34, 9).Value = 5683.5
$ws.Cells.Item(34, 10).Value = 4213.5454
$ws.Cells.Item(34, 11).Value = 5683.5
$ws.Cells.Item(34, 12).Value = 4213.5454
$ws.Cells.Item(34, 13).Value = -5481.5
$ws.Cells.Item(34, 14).Value = -4617.5454

$ws.Cells.Item(58, 8).Value = 1282278
$ws.Cells.Item(58, 9).Value = 1856541.8
$ws.Cells.Item(58, 10).Value = 3236
$ws.Cells.Item(58, 11).Value = 1856541.8
$ws.Cells.Item(58, 12).Value = 3236
$ws.Cells.Item(58, 13).Value = -1856338.8
$ws.Cells.Item(58, 14).Value = -3642

$ws.Cells.Item(127, 8).Value = 79800
$ws.Cells.Item(127, 10).Value = 79800
$ws.Cells.Item(127, 12).Value = 79800
$ws.Cells.Item(127, 14).Value = -89720

$ws.Cells.Item(132, 8).Value = 5488.343
$ws.Cells.Item(132, 9).Value = 7215.4736
$ws.Cells.Item(132, 10).Value = 3437.375
$ws.Cells.Item(132, 11).Value = 21646.4208
$ws.Cells.Item(132, 12).Value = 10312.125
$ws.Cells.Item(132, 13).Value = -19116.4208
$ws.Cells.Item(132, 14).Value = -15372.125

$ws.Cells.Item(134, 8).Value = 3173.3618
$ws.Cells.Item(134, 9).Value = 2244.348
$ws.Cells.Item(134, 10).Value = 4063.6667
$ws.Cells.Item(134, 11).Value = 6733.044
$ws.Cells.Item(134, 12).Value = 12191.0001
$ws.Cells.Item(134, 13).Value = -4198.044
$ws.Cells.Item(134, 14).Value = -17261.0001

$ws.Cells.Item(136, 8).Value = 1282278
$ws.Cells.Item(136, 9).Value = 1856541.8
$ws.Cells.Item(136, 10).Value = 3236
$ws.Cells.Item(136, 11).Value = 5569625.4
$ws.Cells.Item(136, 12).Value = 9708
$ws.Cells.Item(136, 13).Value = -5567075.4
$ws.Cells.Item(136, 14).Value = -14808

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 45024
$ws.Cells.Item(131, 10).Value = 97279.8
$ws.Cells.Item(131, 12).Value = 291839.4
$ws.Cells.Item(131, 14).Value = -301919.4

$ws.Cells.Item(132, 8).Value = 2206.25
$ws.Cells.Item(132, 9).Value = 3200
$ws.Cells.Item(132, 10).Value = 1754.5454
$ws.Cells.Item(132, 11).Value = 28800
$ws.Cells.Item(132, 12).Value = 15790.9086
$ws.Cells.Item(132, 13).Value = -26270
$ws.Cells.Item(132, 14).Value = -20850.9086

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7117.575
$ws.Cells.Item(132, 9).Value = 5738.5
$ws.Cells.Item(132, 11).Value = 17215.5
$ws.Cells.Item(132, 13).Value = -14685.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1530.8636
$ws.Cells.Item(16, 9).Value = 907.13336
$ws.Cells.Item(16, 10).Value = 2867.4285
$ws.Cells.Item(16, 11).Value = 907.13336
$ws.Cells.Item(16, 12).Value = 2867.4285
$ws.Cells.Item(16, 13).Value = -737.13336
$ws.Cells.Item(16, 14).Value = -3207.4285

$ws.Cells.Item(18, 8).Value = 7400
$ws.Cells.Item(18, 9).Value = 5000
$ws.Cells.Item(18, 11).Value = 5000
$ws.Cells.Item(18, 13).Value = -4828

$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 932.1667
$ws.Cells.Item(46, 9).Value = 533.3333
$ws.Cells.Item(46, 10).Value = 1065.1111
$ws.Cells.Item(46, 11).Value = 533.3333
$ws.Cells.Item(46, 12).Value = 1065.1111
$ws.Cells.Item(46, 13).Value = -345.3333
$ws.Cells.Item(46, 14).Value = -1441.1111

$ws.Cells.Item(132, 8).Value = 4511.9556
$ws.Cells.Item(132, 9).Value = 4779.1665
$ws.Cells.Item(132, 10).Value = 3977.5334
$ws.Cells.Item(132, 11).Value = 14337.4995
$ws.Cells.Item(132, 12).Value = 11932.6002
$ws.Cells.Item(132, 13).Value = -11807.4995
$ws.Cells.Item(132, 14).Value = -16992.6002

$ws.Cells.Item(136, 8).Value = 3261.7568
$ws.Cells.Item(136, 9).Value = 1961.0377
$ws.Cells.Item(136, 11).Value = 5883.1131
$ws.Cells.Item(136, 13).Value = -3333.1131

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3354.9048
$ws.Cells.Item(62, 9).Value = 3187.5
$ws.Cells.Item(62, 10).Value = 3457.923
$ws.Cells.Item(62, 11).Value = 3187.5
$ws.Cells.Item(62, 12).Value = 3457.923
$ws.Cells.Item(62, 13).Value = -2563.5
$ws.Cells.Item(62, 14).Value = -4705.923

$ws.Cells.Item(65, 8).Value = 3354.9048
$ws.Cells.Item(65, 9).Value = 3187.5
$ws.Cells.Item(65, 10).Value = 3457.923
$ws.Cells.Item(65, 11).Value = 15937.5
$ws.Cells.Item(65, 12).Value = 17289.615
$ws.Cells.Item(65, 13).Value = -12817.5
$ws.Cells.Item(65, 14).Value = -23529.615

$ws.Cells.Item(132, 8).Value = 1176.1765
$ws.Cells.Item(132, 9).Value = 475.04443
$ws.Cells.Item(132, 10).Value = 2547.9565
$ws.Cells.Item(132, 11).Value = 1425.13329
$ws.Cells.Item(132, 12).Value = 7643.869499999999
$ws.Cells.Item(132, 13).Value = 1104.86671
$ws.Cells.Item(132, 14).Value = -12703.8695

$ws.Cells.Item(136, 8).Value = 4343.5977
$ws.Cells.Item(136, 9).Value = 3290.0393
$ws.Cells.Item(136, 10).Value = 6076.871
$ws.Cells.Item(136, 11).Value = 9870.117899999999
$ws.Cells.Item(136, 12).Value = 18230.613
$ws.Cells.Item(136, 13).Value = -7320.117899999999
$ws.Cells.Item(136, 14).Value = -23330.613
